$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to reflect its new purpose
$ws.Name = "Text_Formats"

# Header row + data rows. Column order on the sheet is TEXT, CAMEL, PASCAL, TITLE.
$data = @(
    @("TEXT",        "CAMEL",      "PASCAL",      "TITLE"),
    @("word",         "word",       "Word",        "Word"),
    @("Title Case",   "titleCase",  "TitleCase",   "Title Case"),
    @("miXed caSE",   "miXedCaSE",  "MiXedCaSE",   "Mixed Case"),
    @("camelCase",    "camelCase",  "CamelCase",   "Camelcase"),
    @("PascalCase",   "pascalCase", "PascalCase",  "Pascalcase"),
    @("lower case",   "lowerCase",  "LowerCase",   "Lower Case"),
    @("UPPER CASE",   "uPPERCASE",  "UPPERCASE",   "Upper Case"),
    @("l",            "l",          "L",           "L")
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}

# Fit the columns to their content
$ws.Columns.Item("A:D").AutoFit() | Out-Null

# Turn the data into a table (ListObject). The original table range included one
# extra (blank) row below the data.
$range = $ws.Range("A1:D10")
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $range, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"
$tbl.TableStyle = "TableStyleLight9"

# Restore the selection that was active when the workbook was last saved
$ws.Range("H9").Select()
